# Update cryptocurrency price/volume figures per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'39.967.23"
$ws.Range("E2").Value = '  -4.25%  '
$ws.Range("D3").Value = "'2.329.48"
$ws.Range("E3").Value = '  -5.96%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = "'306.72"
$ws.Range("D6").Value = "'84.22"
$ws.Range("E6").Value = '  -8.74%  '
$ws.Range("D7").Value = "'0.528"
$ws.Range("E7").Value = '  -4.10%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = "'0.483"
$ws.Range("E9").Value = '  -5.51%  '
$ws.Range("D10").Value = "'0.0822"
$ws.Range("E10").Value = '  -4.70%  '
$ws.Range("D11").Value = "'30.08"
$ws.Range("E11").Value = '  -8.99%  '
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("D13").Value = "'2.686.57"
$ws.Range("E13").Value = '  -6.03%  '
$ws.Range("D14").Value = "'6.40"
$ws.Range("E14").Value = '  -7.13%  '
$ws.Range("D15").Value = "'14.71"
$ws.Range("E15").Value = '  -5.18%  '
$ws.Range("D16").Value = "'2.318.53"
$ws.Range("E16").Value = '  -6.20%  '
$ws.Range("D17").Value = "'0.747"
$ws.Range("E17").Value = '  -6.07%  '
$ws.Range("D18").Value = "'39.907.40"
$ws.Range("E18").Value = '  -4.23%  '
$ws.Range("D19").Value = "'0.0₃0902"
$ws.Range("E19").Value = '  -4.35%  '
$ws.Range("D20").Value = "'6.07"
$ws.Range("E20").Value = '  -5.74%  '
$ws.Range("D21").Value = "'67.60"
$ws.Range("E21").Value = '  -4.31%  '
$ws.Range("D22").Value = "'10.60"
$ws.Range("E22").Value = '  -6.01%  '
$ws.Range("D23").Value = "'235.04"
$ws.Range("E23").Value = '  -1.99%  '
$ws.Range("D24").Value = "'2.54"
$ws.Range("E24").Value = '  -7.80%  '
$ws.Range("E25").Value = '  +0.22%  '
$ws.Range("E26").Value = '  -7.80%  '
$ws.Range("D27").Value = "'23.36"
$ws.Range("E27").Value = '  -6.42%  '
$ws.Range("D28").Value = "'2.14"
$ws.Range("E28").Value = '  -4.36%  '
$ws.Range("D29").Value = "'9.18"
$ws.Range("E29").Value = '  -5.91%  '
$ws.Range("D30").Value = "'34.68"
$ws.Range("E30").Value = '  -5.66%  '
$ws.Range("D31").Value = "'150.86"
$ws.Range("E31").Value = '  -4.20%  '
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("D33").Value = "'5.11"
$ws.Range("E33").Value = '  -6.25%  '
$ws.Range("D34").Value = "'0.0723"
$ws.Range("E34").Value = '  -5.53%  '
$ws.Range("E35").Value = '  -5.65%  '
$ws.Range("E36").Value = '  -2.48%  '
$ws.Range("D37").Value = "'2.77"
$ws.Range("E37").Value = '  -4.05%  '
$ws.Range("D38").Value = "'0.0993"
$ws.Range("E38").Value = '  -4.35%  '
$ws.Range("D39").Value = "'15.70"
$ws.Range("E39").Value = '  -8.71%  '
$ws.Range("D40").Value = "'1.70"
$ws.Range("E40").Value = '  -7.83%  '
$ws.Range("D41").Value = "'3.79"
$ws.Range("E41").Value = '  -5.70%  '
$ws.Range("D42").Value = "'2.29"
$ws.Range("E42").Value = '  -5.20%  '
$ws.Range("D43").Value = "'1.940.57"
$ws.Range("E43").Value = '  -2.62%  '
$ws.Range("D44").Value = "'0.0265"
$ws.Range("E44").Value = '  -6.64%  '
$ws.Range("D45").Value = "'17.64"
$ws.Range("E45").Value = '  -5.88%  '
$ws.Range("D46").Value = "'9.26"
$ws.Range("E46").Value = '  -2.04%  '
$ws.Range("D47").Value = "'2.66"
$ws.Range("E47").Value = '  -10.72%  '
$ws.Range("D48").Value = "'2.548.37"
$ws.Range("E48").Value = '  -6.88%  '
$ws.Range("D49").Value = "'92.36"
$ws.Range("E49").Value = '  -5.35%  '
$ws.Range("D50").Value = "'70.78"
$ws.Range("E50").Value = '  -6.93%  '
$ws.Range("D51").Value = "'63.16"
$ws.Range("E51").Value = '  -6.42%  '
